# Applies the SFDR framework update:
#   "Yes/No/No Evidence Found" quality bucket is replaced by "Yes/No"
#   for every field's Component column on the "Framework Data Model" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Framework Data Model")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row  # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G: "Component"
    if ($cell.Value2 -eq "Yes/No/No Evidence Found") {
        $cell.Value2 = "Yes/No"
    }
}
